$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '61.667.31'
    'E2' = '  -1.95%  '
    'D3' = '2.895.06'
    'E3' = '  -1.88%  '
    'E4' = '  +0.01%  '
    'D5' = '567.67'
    'E5' = '  -4.09%  '
    'D6' = '143.29'
    'E6' = '  -3.20%  '
    'E7' = '  +0.11%  '
    'E8' = '  -0.43%  '
    'D9' = '2.894.29'
    'E9' = '  -1.89%  '
    'D10' = '6.94'
    'E10' = '  -2.25%  '
    'E11' = '  -2.06%  '
    'E12' = '  -2.17%  '
    'E13' = '  -0.97%  '
    'D14' = '31.79'
    'E14' = '  -2.73%  '
    'E15' = '  -0.40%  '
    'D16' = '3.375.67'
    'D17' = '61.690.89'
    'E17' = '  -1.91%  '
    'D18' = '6.54'
    'E18' = '  -1.82%  '
    'D19' = '2.886.66'
    'E19' = '  -2.12%  '
    'D20' = '431.21'
    'E20' = '  -2.11%  '
    'D21' = '13.04'
    'E21' = '  -3.16%  '
    'E22' = '  -1.92%  '
    'E23' = '  -2.53%  '
    'D24' = '79.15'
    'E24' = '  -2.03%  '
    'D25' = '11.88'
    'E25' = '  +0.57%  '
    'E26' = '  +0.02%  '
    'D27' = '9.94'
    'E27' = '  -11.33%  '
    'E28' = '  -5.45%  '
    'E29' = '  +3.54%  '
    'D30' = '7.02'
    'E30' = '  -3.75%  '
    'E31' = '  -4.12%  '
    'E32' = '  -8.56%  '
    'D33' = '1.00'
    'E33' = '  +0.11%  '
    'E34' = '  -1.69%  '
    'E35' = '  -3.16%  '
    'D36' = '0.958'
    'E36' = '  -3.31%  '
    'E37' = '  -4.26%  '
    'D38' = '48.80'
    'E38' = '  -1.85%  '
    'E40' = '  -8.05%  '
    'E41' = '  -3.37%  '
    'E42' = '  -3.51%  '
    'D43' = '39.70'
    'E43' = '  +0.26%  '
    'E44' = '  -4.11%  '
    'D45' = '2.689.39'
    'E45' = '  -0.69%  '
    'D46' = '132.54'
    'E46' = '  -2.33%  '
    'E47' = '  -0.87%  '
    'D48' = '342.91'
    'E48' = '  -5.01%  '
    'E50' = '  -1.54%  '
    'D51' = '21.54'
    'E51' = '  -5.28%  '
}

foreach ($cell in $updates.Keys) {
    $value = $updates[$cell]
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.NumberFormat = "General"
    $r.Style = "Normal"
}
